$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.595.70"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.375.94"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.32%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.635"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("E12").Value = "  -4.07%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "2.735.98"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "2.371.07"
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "42.723.25"
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.38%  "
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "272.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.82%  "
$ws.Range("E24").Value = "  -7.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.88%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.43%  "
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.01%  "
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0911"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0360"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.30%  "
$ws.Range("E40").Value = "  -3.96%  "
$ws.Range("E41").Value = "  +3.25%  "
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.234"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.78%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +28.98%  "
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("D51").Value = "1.586.96"
$ws.Range("E51").Value = "  +6.22%  "
